$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 211, pushing the existing rows
# 211..302 down to 213..304 (dimension grows from R302 to R304).
$ws.Range("A211:A212").EntireRow.Insert()

# New row 211
$ws.Cells.Item(211, 1).Value = 4
$ws.Cells.Item(211, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(211, 3).Value = "Los Lagos"
$ws.Cells.Item(211, 4).Value = 44609
$ws.Cells.Item(211, 5).Value = 10
$ws.Cells.Item(211, 6).Value = 100114013
$ws.Cells.Item(211, 7).Value = "Zanahoria"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 150
$ws.Cells.Item(211, 11).Value = 12000
$ws.Cells.Item(211, 12).Value = 12000
$ws.Cells.Item(211, 13).Value = 12000
$ws.Cells.Item(211, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(211, 15).Value = "Chillán"
$ws.Cells.Item(211, 16).Value = 600
$ws.Cells.Item(211, 17).Value = 20
$ws.Cells.Item(211, 18).Value = "Hortaliza"

# New row 212
$ws.Cells.Item(212, 1).Value = 4
$ws.Cells.Item(212, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(212, 3).Value = "Los Lagos"
$ws.Cells.Item(212, 4).Value = 44609
$ws.Cells.Item(212, 5).Value = 10
$ws.Cells.Item(212, 6).Value = 100114013
$ws.Cells.Item(212, 7).Value = "Zanahoria"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 150
$ws.Cells.Item(212, 11).Value = 10000
$ws.Cells.Item(212, 12).Value = 10000
$ws.Cells.Item(212, 13).Value = 10000
$ws.Cells.Item(212, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(212, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(212, 16).Value = 500
$ws.Cells.Item(212, 17).Value = 20
$ws.Cells.Item(212, 18).Value = "Hortaliza"
